$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 1432.8334
$ws.Cells.Item(9, 9).Value = 299.33334
$ws.Cells.Item(9, 11).Value = 299.33334
$ws.Cells.Item(9, 13).Value = -130.33334
$ws.Cells.Item(11, 8).Value = 757.375
$ws.Cells.Item(11, 9).Value = 757.375
$ws.Cells.Item(11, 11).Value = 757.375
$ws.Cells.Item(11, 13).Value = -617.375
$ws.Cells.Item(18, 8).Value = 769.8
$ws.Cells.Item(18, 9).Value = 749.6667
$ws.Cells.Item(18, 10).Value = 800
$ws.Cells.Item(18, 11).Value = 749.6667
$ws.Cells.Item(18, 12).Value = 800
$ws.Cells.Item(18, 13).Value = -465.6667
$ws.Cells.Item(18, 14).Value = -1368
$ws.Cells.Item(64, 8).Value = 12999.167
$ws.Cells.Item(64, 10).Value = 19000
$ws.Cells.Item(64, 12).Value = 19000
$ws.Cells.Item(64, 14).Value = -19496
$ws.Cells.Item(67, 8).Value = 12999.167
$ws.Cells.Item(67, 10).Value = 19000
$ws.Cells.Item(67, 12).Value = 19000
$ws.Cells.Item(67, 14).Value = -20716
$ws.Cells.Item(74, 8).Value = 4330
$ws.Cells.Item(74, 9).Value = 4362.5
$ws.Cells.Item(74, 10).Value = 4200
$ws.Cells.Item(74, 11).Value = 4362.5
$ws.Cells.Item(74, 12).Value = 4200
$ws.Cells.Item(74, 13).Value = -3426.5
$ws.Cells.Item(74, 14).Value = -6072
$ws.Cells.Item(77, 8).Value = 4330
$ws.Cells.Item(77, 9).Value = 4362.5
$ws.Cells.Item(77, 10).Value = 4200
$ws.Cells.Item(77, 11).Value = 21812.5
$ws.Cells.Item(77, 12).Value = 21000
$ws.Cells.Item(77, 13).Value = -17132.5
$ws.Cells.Item(77, 14).Value = -30360
$ws.Cells.Item(109, 8).Value = 44936.8
$ws.Cells.Item(109, 10).Value = 44936.8
$ws.Cells.Item(109, 12).Value = 44936.8
$ws.Cells.Item(109, 14).Value = -47710.8
$ws.Cells.Item(116, 8).Value = 6426
$ws.Cells.Item(116, 9).Value = 6634.6665
$ws.Cells.Item(116, 11).Value = 6634.6665
$ws.Cells.Item(116, 13).Value = -3192.6665
$ws.Cells.Item(121, 8).Value = 897
$ws.Cells.Item(121, 10).Value = 897
$ws.Cells.Item(121, 12).Value = 2691
$ws.Cells.Item(121, 14).Value = -6185
$ws.Cells.Item(138, 8).Value = 13359.814
$ws.Cells.Item(138, 10).Value = 13408.27
$ws.Cells.Item(138, 12).Value = 40224.81
$ws.Cells.Item(138, 14).Value = -50504.81
$ws.Cells.Item(141, 8).Value = 1199.6666
$ws.Cells.Item(141, 9).Value = 1199.6666
$ws.Cells.Item(141, 11).Value = 3598.9998
$ws.Cells.Item(141, 13).Value = 1581.0002
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3415.8333
$ws.Cells.Item(61, 9).Value = 1665.3334
$ws.Cells.Item(61, 10).Value = 5166.3335
$ws.Cells.Item(61, 11).Value = 1665.3334
$ws.Cells.Item(61, 12).Value = 5166.3335
$ws.Cells.Item(61, 13).Value = -1453.3334
$ws.Cells.Item(61, 14).Value = -5590.3335
$ws.Cells.Item(88, 8).Value = 6000.8
$ws.Cells.Item(88, 9).Value = 5006
$ws.Cells.Item(88, 10).Value = 6249.5
$ws.Cells.Item(88, 11).Value = 5006
$ws.Cells.Item(88, 12).Value = 6249.5
$ws.Cells.Item(88, 13).Value = -4600
$ws.Cells.Item(88, 14).Value = -7061.5
$ws.Cells.Item(91, 8).Value = 6000.8
$ws.Cells.Item(91, 9).Value = 5006
$ws.Cells.Item(91, 10).Value = 6249.5
$ws.Cells.Item(91, 11).Value = 5006
$ws.Cells.Item(91, 12).Value = 6249.5
$ws.Cells.Item(91, 13).Value = -3602
$ws.Cells.Item(91, 14).Value = -9057.5
$ws.Cells.Item(102, 8).Value = 1483.3334
$ws.Cells.Item(102, 9).Value = 1483.3334
$ws.Cells.Item(102, 11).Value = 1483.3334
$ws.Cells.Item(102, 13).Value = 138.6666
$ws.Cells.Item(132, 8).Value = 1809.8182
$ws.Cells.Item(132, 9).Value = 1410.375
$ws.Cells.Item(132, 11).Value = 4231.125
$ws.Cells.Item(132, 13).Value = -1701.125
$ws.Cells.Item(136, 8).Value = 3415.8333
$ws.Cells.Item(136, 9).Value = 1665.3334
$ws.Cells.Item(136, 10).Value = 5166.3335
$ws.Cells.Item(136, 11).Value = 4996.0002
$ws.Cells.Item(136, 12).Value = 15499.0005
$ws.Cells.Item(136, 13).Value = -2446.0002
$ws.Cells.Item(136, 14).Value = -20599.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1800
$ws.Cells.Item(86, 10).Value = 1800
$ws.Cells.Item(86, 12).Value = 1800
$ws.Cells.Item(86, 14).Value = -4046
$ws.Cells.Item(89, 8).Value = 1800
$ws.Cells.Item(89, 10).Value = 1800
$ws.Cells.Item(89, 12).Value = 9000
$ws.Cells.Item(89, 14).Value = -20232
$ws.Cells.Item(99, 8).Value = 1699.2
$ws.Cells.Item(99, 9).Value = 1570.7142
$ws.Cells.Item(99, 10).Value = 1999
$ws.Cells.Item(99, 11).Value = 1570.7142
$ws.Cells.Item(99, 12).Value = 1999
$ws.Cells.Item(99, 13).Value = -72.71419999999989
$ws.Cells.Item(99, 14).Value = -4995
$ws.Cells.Item(107, 8).Value = 1431.5
$ws.Cells.Item(107, 9).Value = 1431.5
$ws.Cells.Item(107, 11).Value = 1431.5
$ws.Cells.Item(107, 13).Value = 488.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3547.8
$ws.Cells.Item(62, 9).Value = 2446.3333
$ws.Cells.Item(62, 11).Value = 2446.3333
$ws.Cells.Item(62, 13).Value = -1822.3333
$ws.Cells.Item(65, 8).Value = 3547.8
$ws.Cells.Item(65, 9).Value = 2446.3333
$ws.Cells.Item(65, 11).Value = 12231.6665
$ws.Cells.Item(65, 13).Value = -9111.666499999999
$ws.Cells.Item(97, 8).Value = 16000
$ws.Cells.Item(97, 10).Value = 16000
$ws.Cells.Item(97, 12).Value = 16000
$ws.Cells.Item(97, 14).Value = -17982
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 831.2857
$ws.Cells.Item(12, 9).Value = 130
$ws.Cells.Item(12, 10).Value = 1111.8
$ws.Cells.Item(12, 11).Value = 390
$ws.Cells.Item(12, 12).Value = 3335.4
$ws.Cells.Item(12, 13).Value = -217
$ws.Cells.Item(12, 14).Value = -3681.4
$ws.Cells.Item(137, 8).Value = 2000
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 14).Value = -16200
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82, 8).Value = 29333
$ws.Cells.Item(82, 10).Value = 29333
$ws.Cells.Item(82, 12).Value = 29333
$ws.Cells.Item(82, 14).Value = -30099
$ws.Cells.Item(85, 8).Value = 29333
$ws.Cells.Item(85, 10).Value = 29333
$ws.Cells.Item(85, 12).Value = 29333
$ws.Cells.Item(85, 14).Value = -31985
$ws.Cells.Item(100, 8).Value = 6971595
$ws.Cells.Item(100, 9).Value = 9957722
$ws.Cells.Item(100, 11).Value = 19915444
$ws.Cells.Item(100, 13).Value = -19914903
$ws.Cells.Item(122, 8).Value = 901
$ws.Cells.Item(122, 9).Value = 602
$ws.Cells.Item(122, 11).Value = 1806
$ws.Cells.Item(122, 13).Value = 644
$ws.Cells.Item(132, 8).Value = 3056.182
$ws.Cells.Item(132, 9).Value = 1916.625
$ws.Cells.Item(132, 11).Value = 5749.875
$ws.Cells.Item(132, 13).Value = -3219.875
$ws.Cells.Item(136, 8).Value = 45415.832
$ws.Cells.Item(136, 9).Value = 58332.332
$ws.Cells.Item(136, 11).Value = 174996.996
$ws.Cells.Item(136, 13).Value = -172446.996

Write-Output "Applied 159 cell updates"